$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.936.76'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '''1.671.27'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''214.76'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').Value = '''20.19'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '''1.906.90'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '''1.670.15'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').Value = '''0.526'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').Value = '''65.50'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '''26.927.44'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '''8.05'
$ws.Range('E18').Value = '  +3.97%  '
$ws.Range('D19').Value = '''233.53'
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D22').Value = '''4.42'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').Value = '''9.17'
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').Value = '''145.83'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').Value = '''7.12'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').Value = '''15.92'
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').Value = '''3.33'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').Value = '''1.459.94'
$ws.Range('E33').Value = '  -5.64%  '
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('D35').Value = '''1.65'
$ws.Range('E35').Value = '  +1.88%  '
$ws.Range('D36').Value = '''2.42'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '''0.580'
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('D38').Value = '''0.899'
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('D39').Value = '''0.0171'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('E40').Value = '  +13.48%  '
$ws.Range('D41').Value = '''5.77'
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('D44').Value = '''66.27'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('D45').Value = '''1.811.50'
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('D46').Value = '''0.781'
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('D47').Value = '''90.61'
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('D48').Value = '''1.54'
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').Value = '''0.102'
$ws.Range('E49').Value = '  +2.79%  '
$ws.Range('D50').Value = '''0.0508'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').Value = '''7.63'
$ws.Range('E51').Value = '  -0.38%  '
